$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# auto-converted to numbers by Excel, then restore default styling so no
# stray style index is left applied to the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.026.70"
$ws.Range("E2").Value = "  -6.73%  "
$ws.Range("D3").Value = "3.490.37"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "389.97"
$ws.Range("E5").Value = "  -6.57%  "
$ws.Range("D6").Value = "119.62"
$ws.Range("E6").Value = "  -8.15%  "
$ws.Range("D7").Value = "3.480.65"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -9.91%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "0.670"
$ws.Range("E10").Value = "  -12.58%  "
$ws.Range("D11").Value = "0.148"
$ws.Range("E11").Value = "  -16.33%  "
$ws.Range("D12").Value = "0.0000324"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("D13").Value = "38.51"
$ws.Range("E13").Value = "  -9.08%  "
$ws.Range("D14").Value = "4.063.99"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "9.11"
$ws.Range("E15").Value = "  -8.92%  "
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "3.486.69"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "18.56"
$ws.Range("E19").Value = "  -9.22%  "
$ws.Range("D20").Value = "63.124.00"
$ws.Range("E20").Value = "  -6.43%  "
$ws.Range("D21").Value = "1.01"
$ws.Range("E21").Value = "  -12.00%  "
$ws.Range("D22").Value = "390.20"
$ws.Range("E22").Value = "  -15.36%  "
$ws.Range("D23").Value = "13.78"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").Value = "80.50"
$ws.Range("E24").Value = "  -8.65%  "
$ws.Range("D25").Value = "2.83"
$ws.Range("E25").Value = "  -9.08%  "
$ws.Range("D26").Value = "32.97"
$ws.Range("E26").Value = "  -5.84%  "
$ws.Range("D27").Value = "5.14"
$ws.Range("E27").Value = "  +5.72%  "
$ws.Range("D28").Value = "2.96"
$ws.Range("E28").Value = "  -11.86%  "
$ws.Range("D29").Value = "8.69"
$ws.Range("E29").Value = "  -14.82%  "
$ws.Range("D30").Value = "11.79"
$ws.Range("E30").Value = "  -5.01%  "
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").Value = "  -7.39%  "
$ws.Range("D32").Value = "0.109"
$ws.Range("E32").Value = "  -6.88%  "
$ws.Range("D33").Value = "6.78"
$ws.Range("E33").Value = "  -8.69%  "
$ws.Range("D34").Value = "0.149"
$ws.Range("E34").Value = "  -7.94%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("D36").Value = "36.35"
$ws.Range("E36").Value = "  -13.13%  "
$ws.Range("D37").Value = "53.56"
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("D38").Value = "0.0433"
$ws.Range("E38").Value = "  -12.23%  "
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  +14.44%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.130"
$ws.Range("E41").Value = "  -11.60%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").Value = "  +13.53%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0618"
$ws.Range("E43").Value = "  -13.21%  "
$ws.Range("D44").Value = "140.88"
$ws.Range("E44").Value = "  -5.33%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  -8.92%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "2.70"
$ws.Range("E46").Value = "  -11.27%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "24.62"
$ws.Range("E47").Value = "  +14.08%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "3.05"
$ws.Range("E48").Value = "  -6.60%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "3.98"
$ws.Range("E50").Value = "  -7.79%  "
$ws.Range("D51").Value = "0.274"
$ws.Range("E51").Value = "  -11.78%  "

# Remove the temporary number-format override so the cells fall back to the
# workbook's default (unstyled) cell format, matching the original layout.
$priceRange.Style = "Normal"

